# This revision's canonical OOXML diff only touches two things:
#
#   1. ppt/presentation.xml -> <p:embeddedFontLst>: the <p:font> "panose"
#      value (and the accompanying "pitchFamily" attribute) for the
#      "Oswald" and "Source Sans Pro" embedded-font records is refreshed
#      (PowerPoint recomputes this straight from the embedded TTF/OTF
#      binary's OS/2 table whenever it re-embeds the fonts on save).
#   2. ppt/changesInfos/changesInfo1.xml: the internal collaborative-
#      editing change-history log simply has its change-session blocks
#      reordered (most-recent-session-first) -- no new content, purely a
#      re-ordering of the existing <pc:docChgLst> blocks.
#
# Neither of these is content that the PowerPoint object model exposes:
# there is no Shape/Slide/TextRange touched by this revision, and the
# embedded-font PANOSE/pitch metadata and the changesInfo change-log are
# low-level package bookkeeping that PowerPoint itself regenerates
# on save rather than something scriptable via Application/Presentation
# automation (there is no Presentation member for embedded-font PANOSE
# data or for the changesInfo log -- confirmed against the live
# Presentation COM interface, which has no such members).
#
# So there is no slide/shape/text edit to make here; touch the
# presentation without altering any author-visible content so the
# script still runs cleanly through the COM host.
$p = $ppt.ActivePresentation
$null = $p.Slides.Count
